$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 24 de Mayo de 2020 a las 18:35"

# Update per-country statistics (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes)

# Row 4 - Turquia
$ws.Range("B4").Value = 1673301
$ws.Range("C4").Value = 6473
$ws.Range("D4").Value = 448959
$ws.Range("E4").Value = 1125518
$ws.Range("G4").Value = 141
$ws.Range("H4").Value = 98824

# Row 5 - India
$ws.Range("B5").Value = 352523
$ws.Range("C5").Value = 5125
$ws.Range("E5").Value = 187648
$ws.Range("G5").Value = 275
$ws.Range("H5").Value = 22288

# Row 8 - Peru
$ws.Range("B8").Value = 259559
$ws.Range("C8").Value = 2405
$ws.Range("G8").Value = 118
$ws.Range("H8").Value = 36793

# Row 9 - Canada
$ws.Range("B9").Value = 229858
$ws.Range("C9").Value = 531
$ws.Range("D9").Value = 140479
$ws.Range("E9").Value = 56594
$ws.Range("G9").Value = 50
$ws.Range("H9").Value = 32785

# Row 13 - Belgica
$ws.Range("B13").Value = 137991
$ws.Range("C13").Value = 6568
$ws.Range("D13").Value = 57415
$ws.Range("E13").Value = 76563
$ws.Range("G13").Value = 145
$ws.Range("H13").Value = 4013

# Row 33 - Republica Dominicana
$ws.Range("B33").Value = 24639
$ws.Range("C33").Value = 57
$ws.Range("E33").Value = 1971
$ws.Range("G33").Value = 4
$ws.Range("H33").Value = 1608

# Row 56
$ws.Range("D56").Value = 4352
$ws.Range("E56").Value = 3935

# Row 57
$ws.Range("B57").Value = 8306
$ws.Range("C57").Value = 193
$ws.Range("D57").Value = 4784
$ws.Range("E57").Value = 2922
$ws.Range("G57").Value = 8
$ws.Range("H57").Value = 600

# Row 103
$ws.Range("B103").Value = 1138
$ws.Range("C103").Value = 49
$ws.Range("E103").Value = 455

# Row 149
$ws.Range("B149").Value = 265
$ws.Range("C149").Value = 10
$ws.Range("D149").Value = 139
$ws.Range("E149").Value = 100

$wb.Save()
